$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated values (new TPM-derived figures) for rows 2-6, columns M-T
$updates = @{
    "M2" = 0.939124
    "N2" = 2.817372
    "O2" = 0.4188652080327055
    "P2" = 0.4188652080327054
    "Q2" = 0.05168781975333333
    "R2" = 0.46519037778
    "S2" = 0.4188652080327055
    "T2" = 0.4188652080327054

    "O3" = 0.3168364554023206
    "P3" = 0.3168364554023206
    "S3" = 0.3168364554023206
    "T3" = 0.3168364554023206

    "M4" = 0.2295973333333333
    "N4" = 0.6887920000000001
    "O4" = 0.102404298889626
    "P4" = 0.102404298889626
    "Q4" = 0.01263665456444444
    "R4" = 0.11372989108
    "S4" = 0.102404298889626
    "T4" = 0.102404298889626

    "M5" = 0.2429623333333333
    "N5" = 0.728887
    "O5" = 0.1083653152254422
    "P5" = 0.1083653152254422
    "Q5" = 0.01337224188944444
    "R5" = 0.120350177005
    "S5" = 0.1083653152254422
    "T5" = 0.1083653152254422

    "M6" = 0.120015
    "N6" = 0.3600450000000001
    "O6" = 0.05352872244990561
    "P6" = 0.05352872244990561
    "Q6" = 0.006605425575000001
    "R6" = 0.05944883017500002
    "S6" = 0.05352872244990561
    "T6" = 0.05352872244990561
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
